$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H33").Value = 209.76471
$ws.Range("I33").Value = 189.75
$ws.Range("K33").Value = 189.75
$ws.Range("M33").Value = 39.25
$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5877
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -29384
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 644.2857
$ws.Range("I107").Value = 697.43475
$ws.Range("J107").Value = 399.8
$ws.Range("K107").Value = 697.43475
$ws.Range("L107").Value = 399.8
$ws.Range("M107").Value = 1222.56525
$ws.Range("N107").Value = -4239.8
$ws.Range("H132").Value = 1606.2025
$ws.Range("I132").Value = 1597.1487
$ws.Range("J132").Value = 1740.2
$ws.Range("K132").Value = 4791.4461
$ws.Range("L132").Value = 5220.6
$ws.Range("M132").Value = -2261.4461
$ws.Range("N132").Value = -10280.6
$ws.Range("H135").Value = 1322.2162
$ws.Range("I135").Value = 1301.8214
$ws.Range("K135").Value = 11716.3926
$ws.Range("M135").Value = -9181.392600000001
$ws.Range("H137").Value = 2255.0952
$ws.Range("I137").Value = 2123.2307
$ws.Range("J137").Value = 2469.375
$ws.Range("K137").Value = 6369.6921
$ws.Range("L137").Value = 7408.125
$ws.Range("M137").Value = -3819.6921
$ws.Range("N137").Value = -12508.125
$ws.Range("H138").Value = 3887.7778
$ws.Range("I138").Value = 2052.7368
$ws.Range("J138").Value = 4323.6
$ws.Range("K138").Value = 6158.2104
$ws.Range("L138").Value = 12970.8
$ws.Range("M138").Value = -1018.2104
$ws.Range("N138").Value = -23250.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13529.906
$ws.Range("I32").Value = 7010.5713
$ws.Range("J32").Value = 47756.418
$ws.Range("K32").Value = 7010.5713
$ws.Range("L32").Value = 47756.418
$ws.Range("M32").Value = -6723.5713
$ws.Range("N32").Value = -48330.418
$ws.Range("H37").Value = 19857
$ws.Range("H45").Value = 9592.462
$ws.Range("I45").Value = 14156.875
$ws.Range("K45").Value = 14156.875
$ws.Range("M45").Value = -13779.875
$ws.Range("H61").Value = 300054.62
$ws.Range("I61").Value = 2975.9546
$ws.Range("K61").Value = 2975.9546
$ws.Range("M61").Value = -2763.9546
$ws.Range("H74").Value = 14989.9375
$ws.Range("I74").Value = 3737
$ws.Range("J74").Value = 48748.75
$ws.Range("K74").Value = 3737
$ws.Range("L74").Value = 48748.75
$ws.Range("M74").Value = -2863
$ws.Range("N74").Value = -50496.75
$ws.Range("H77").Value = 14989.9375
$ws.Range("I77").Value = 3737
$ws.Range("J77").Value = 48748.75
$ws.Range("K77").Value = 18685
$ws.Range("L77").Value = 243743.75
$ws.Range("M77").Value = -14317
$ws.Range("N77").Value = -252479.75
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1580.4286
$ws.Range("I122").Value = 1409.45
$ws.Range("K122").Value = 4228.35
$ws.Range("M122").Value = -1778.35
$ws.Range("H132").Value = 2522.5227
$ws.Range("I132").Value = 2108.303
$ws.Range("K132").Value = 6324.909
$ws.Range("M132").Value = -3794.909
$ws.Range("H136").Value = 300054.62
$ws.Range("I136").Value = 2975.9546
$ws.Range("K136").Value = 8927.863799999999
$ws.Range("M136").Value = -6377.863799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3779.875
$ws.Range("I20").Value = 4325.4287
$ws.Range("K20").Value = 4325.4287
$ws.Range("M20").Value = -4078.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29059.736
$ws.Range("I31").Value = 35368.066
$ws.Range("J31").Value = 5403.5
$ws.Range("K31").Value = 35368.066
$ws.Range("L31").Value = 5403.5
$ws.Range("M31").Value = -35073.066
$ws.Range("N31").Value = -5993.5
$ws.Range("H34").Value = 29059.736
$ws.Range("I34").Value = 35368.066
$ws.Range("J34").Value = 5403.5
$ws.Range("K34").Value = 35368.066
$ws.Range("L34").Value = 5403.5
$ws.Range("M34").Value = -35166.066
$ws.Range("N34").Value = -5807.5
$ws.Range("H58").Value = 2015.2812
$ws.Range("I58").Value = 1664.3214
$ws.Range("K58").Value = 1664.3214
$ws.Range("M58").Value = -1461.3214
$ws.Range("H136").Value = 2015.2812
$ws.Range("I136").Value = 1664.3214
$ws.Range("K136").Value = 4992.9642
$ws.Range("M136").Value = -2442.9642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8475
$ws.Range("I3").Value = 1950
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 5850
$ws.Range("L3").Value = 45000
$ws.Range("M3").Value = -5738
$ws.Range("N3").Value = -45224
$ws.Range("H29").Value = 809.46155
$ws.Range("J29").Value = 1125
$ws.Range("L29").Value = 3375
$ws.Range("N29").Value = -3929
$ws.Range("H39").Value = 16832.916
$ws.Range("J39").Value = 16832.916
$ws.Range("L39").Value = 50498.74800000001
$ws.Range("N39").Value = -51086.74800000001
$ws.Range("H68").Value = 16666934
$ws.Range("I68").Value = 401
$ws.Range("K68").Value = 1203
$ws.Range("M68").Value = -392
$ws.Range("H71").Value = 16666934
$ws.Range("I71").Value = 401
$ws.Range("K71").Value = 3609
$ws.Range("M71").Value = 447
$ws.Range("H107").Value = 284.95
$ws.Range("J107").Value = 201.3
$ws.Range("L107").Value = 603.9000000000001
$ws.Range("N107").Value = -4443.9
$ws.Range("H112").Value = 5493.75
$ws.Range("I112").Value = 5493.75
$ws.Range("K112").Value = 16481.25
$ws.Range("M112").Value = -15373.25
$ws.Range("H121").Value = 189049.5
$ws.Range("J121").Value = 302199.6
$ws.Range("L121").Value = 906598.7999999999
$ws.Range("N121").Value = -909218.7999999999
$ws.Range("H123").Value = 895
$ws.Range("I123").Value = 895
$ws.Range("K123").Value = 2685
$ws.Range("M123").Value = -235
$ws.Range("H131").Value = 43230.81
$ws.Range("J131").Value = 4960.04
$ws.Range("L131").Value = 14880.12
$ws.Range("N131").Value = -24960.12
$ws.Range("H132").Value = 1551.3438
$ws.Range("I132").Value = 1399.8182
$ws.Range("J132").Value = 1884.7
$ws.Range("K132").Value = 12598.3638
$ws.Range("L132").Value = 16962.3
$ws.Range("M132").Value = -10068.3638
$ws.Range("N132").Value = -22022.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 99999.586
$ws.Range("I64").Value = 99999.586
$ws.Range("K64").Value = 99999.586
$ws.Range("M64").Value = -99751.586
$ws.Range("H67").Value = 99999.586
$ws.Range("I67").Value = 99999.586
$ws.Range("K67").Value = 99999.586
$ws.Range("M67").Value = -99141.586
$ws.Range("H80").Value = 4511.8887
$ws.Range("I80").Value = 3249.25
$ws.Range("J80").Value = 5522
$ws.Range("K80").Value = 3249.25
$ws.Range("L80").Value = 5522
$ws.Range("M80").Value = -2251.25
$ws.Range("N80").Value = -7518
$ws.Range("H83").Value = 4511.8887
$ws.Range("I83").Value = 3249.25
$ws.Range("J83").Value = 5522
$ws.Range("K83").Value = 16246.25
$ws.Range("L83").Value = 27610
$ws.Range("M83").Value = -11254.25
$ws.Range("N83").Value = -37594
$ws.Range("H132").Value = 2107.0264
$ws.Range("I132").Value = 1469.8182
$ws.Range("J132").Value = 2983.1875
$ws.Range("K132").Value = 4409.4546
$ws.Range("L132").Value = 8949.5625
$ws.Range("M132").Value = -1879.4546
$ws.Range("N132").Value = -14009.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1894.3077
$ws.Range("I93").Value = 1653.6666
$ws.Range("J93").Value = 2100.5715
$ws.Range("K93").Value = 1653.6666
$ws.Range("L93").Value = 2100.5715
$ws.Range("M93").Value = -405.6666
$ws.Range("N93").Value = -4596.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2322.7273
$ws.Range("J96").Value = 2647.25
$ws.Range("L96").Value = 2647.25
$ws.Range("N96").Value = -5393.25
